$d = $word.ActiveDocument

# Remove the "/SKD/TCF/VII/${year}" suffix that used to follow "${no_surat}",
# leaving only the "${no_surat}" placeholder in that paragraph.
$range = $d.Content
$range.Find.Execute("/SKD/TCF/VII/`${year}", $true, $false, $false, $false, $false, `
                     $true, 1, $false, "", 2)
